# Update attendance ("想去人数") figures on the 展览 and 全部类型 sheets,
# plus the lowest ticket price ("最低票价") for row 2 on 全部类型.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1): update column F only ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 2196
$ws1.Range("F3").Value = 628
$ws1.Range("F4").Value = 1578
$ws1.Range("F5").Value = 7366
$ws1.Range("F7").Value = 184

# --- Sheet "全部类型" (sheet4): update column F, and G2 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 2196
$ws4.Range("G2").Value = 60
$ws4.Range("F3").Value = 628
$ws4.Range("F4").Value = 1578
$ws4.Range("F5").Value = 7366
$ws4.Range("F7").Value = 184
